$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "VENTAS POR GRUPO" (sheet1)
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")

# Insert a brand-new data row for "MUÑOZ CALDERON JUAN ADOLFO" just above
# the existing "PALMA PICO OSCAR FILIDEL" row (row 317), shifting every
# row below it down by one.
$ws1.Rows.Item(317).Insert()
$ws1.Cells.Item(317, 1).Value = "OFICINA-CATAECSA"
$ws1.Cells.Item(317, 2).Value = "MUÑOZ CALDERON JUAN ADOLFO"
for ($c = 3; $c -le 18; $c++) {
    $ws1.Cells.Item(317, $c).Value = 0
}

# Individual value corrections elsewhere in the sheet (rows unaffected by
# the insertion above since they sit above row 317).
$ws1.Cells.Item(63, 4).Value = 2136.96    # D63
$ws1.Cells.Item(63, 12).Value = 1108.6    # L63
$ws1.Cells.Item(85, 13).Value = 4804.82   # M85

# The trailing "N de 353" tally row (previously row 355) is now row 356;
# the denominator grows to 354 (one more data row) and the D/L numerators
# each grow by one because D63 and L63 went from 0 to non-zero.
$ws1.Cells.Item(356, 3).Value = "3 de 354"
$ws1.Cells.Item(356, 4).Value = "15 de 354"
$ws1.Cells.Item(356, 5).Value = "6 de 354"
$ws1.Cells.Item(356, 6).Value = "0 de 354"
$ws1.Cells.Item(356, 7).Value = "0 de 354"
$ws1.Cells.Item(356, 8).Value = "5 de 354"
$ws1.Cells.Item(356, 9).Value = "10 de 354"
$ws1.Cells.Item(356, 10).Value = "0 de 354"
$ws1.Cells.Item(356, 11).Value = "3 de 354"
$ws1.Cells.Item(356, 12).Value = "23 de 354"
$ws1.Cells.Item(356, 13).Value = "44 de 354"
$ws1.Cells.Item(356, 14).Value = "2 de 354"
$ws1.Cells.Item(356, 15).Value = "0 de 354"
$ws1.Cells.Item(356, 16).Value = "1 de 354"
$ws1.Cells.Item(356, 17).Value = "0 de 354"
$ws1.Cells.Item(356, 18).Value = "0 de 354"

# ---------------------------------------------------------------------
# Sheet "VENTA MENSUAL" (sheet2)
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")

# Same new-client row, inserted above the existing "PALMA PICO OSCAR
# FILIDEL" row (row 321 in this sheet).
$ws2.Rows.Item(321).Insert()
$ws2.Cells.Item(321, 1).Value = "OFICINA-CATAECSA"
$ws2.Cells.Item(321, 2).Value = "MUÑOZ CALDERON JUAN ADOLFO"
for ($c = 3; $c -le 7; $c++) {
    $ws2.Cells.Item(321, $c).Value = 0
}

# Individual value corrections (rows above the insertion point).
$ws2.Cells.Item(63, 6).Value = 3314.82    # F63 (noviembre)
$ws2.Cells.Item(85, 6).Value = 4804.82    # F85 (noviembre)

# Grand-total row (previously row 359, now row 360): only the "noviembre"
# column total moves, by the same amount the two corrections above added.
$ws2.Cells.Item(360, 6).Value = 143309.12

# ---------------------------------------------------------------------
# Sheet "CUMPLIMIENTO MENSUAL" (sheet3) - derived roll-up values
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# Row 15: CASTRO ALCIVAR EDA MARIA / 240X80 PORCELANATO
$ws3.Cells.Item(15, 4).Value = 7973.56
$ws3.Cells.Item(15, 5).Value = 6851.849999999999
$ws3.Cells.Item(15, 6).Value = 0.5378306569599087

# Row 23: CASTRO ALCIVAR EDA MARIA / PIEDRA SINTERIZADA
$ws3.Cells.Item(23, 4).Value = 5020.76
$ws3.Cells.Item(23, 5).Value = 11127.24
$ws3.Cells.Item(23, 6).Value = 0.3109214763438197

# Row 24: CASTRO ALCIVAR EDA MARIA / PORCELANATO
$ws3.Cells.Item(24, 4).Value = 11135.11
$ws3.Cells.Item(24, 5).Value = 39171.89
$ws3.Cells.Item(24, 6).Value = 0.2213431530403324

# Row 77: TOTAL
$ws3.Cells.Item(77, 4).Value = 143174.66
$ws3.Cells.Item(77, 5).Value = 274074.0197415455
$ws3.Cells.Item(77, 6).Value = 0.3431398754543359
